$d = $word.ActiveDocument

$pairs = @(
    @("29×79=", "55×64="),
    @("21×30=", "39×60="),
    @("74×71=", "13×29="),
    @("80×73=", "85×75="),
    @("11×22=", "89×41="),
    @("59×18=", "87×74="),
    @("45×84=", "51×33="),
    @("27×25=", "69×29="),
    @("13×96=", "33×96="),
    @("44×17=", "56×25="),
    @("38×76=", "70×28="),
    @("43×91=", "70×49="),
    @("21×14=", "71×94="),
    @("54×97=", "56×57="),
    @("73×75=", "41×39="),
    @("59×87=", "74×28="),
    @("20×67=", "65×99="),
    @("61×67=", "31×28="),
    @("37×22=", "55×72="),
    @("26×36=", "62×15="),
    @("39×89=", "83×71="),
    @("16×87=", "21×32="),
    @("27×90=", "24×44="),
    @("59×28=", "11×62="),
    @("84×70=", "64×70=")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
